$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2533.1538
$ws.Range("I40").Value = 5259.6665
$ws.Range("J40").Value = 1715.2
$ws.Range("K40").Value = 5259.6665
$ws.Range("L40").Value = 1715.2
$ws.Range("M40").Value = -5084.6665
$ws.Range("N40").Value = -2065.2
# Row 64
$ws.Range("H64").Value = 2756.25
$ws.Range("I64").Value = 2500
$ws.Range("J64").Value = 2815.3845
$ws.Range("K64").Value = 2500
$ws.Range("L64").Value = 2815.3845
$ws.Range("M64").Value = -2252
$ws.Range("N64").Value = -3311.3845
# Row 67
$ws.Range("H67").Value = 2756.25
$ws.Range("I67").Value = 2500
$ws.Range("J67").Value = 2815.3845
$ws.Range("K67").Value = 2500
$ws.Range("L67").Value = 2815.3845
$ws.Range("M67").Value = -1642
$ws.Range("N67").Value = -4531.3845
# Row 76
$ws.Range("H76").Value = 46139.26
$ws.Range("I76").Value = 65062.688
$ws.Range("J76").Value = 2885.7144
$ws.Range("K76").Value = 65062.688
$ws.Range("L76").Value = 2885.7144
$ws.Range("M76").Value = -64747.688
$ws.Range("N76").Value = -3515.7144
# Row 79
$ws.Range("H79").Value = 46139.26
$ws.Range("I79").Value = 65062.688
$ws.Range("J79").Value = 2885.7144
$ws.Range("K79").Value = 65062.688
$ws.Range("L79").Value = 2885.7144
$ws.Range("M79").Value = -63970.688
$ws.Range("N79").Value = -5069.7144
# Row 132
$ws.Range("H132").Value = 2269271.8
$ws.Range("I132").Value = 2599007.5
$ws.Range("K132").Value = 7797022.5
$ws.Range("M132").Value = -7794492.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 17544912
$ws.Range("I45").Value = 41667340
$ws.Range("J45").Value = 1327.2727
$ws.Range("K45").Value = 41667340
$ws.Range("L45").Value = 1327.2727
$ws.Range("M45").Value = -41666963
$ws.Range("N45").Value = -2081.2727
# Row 63
$ws.Range("H63").Value = 2500850
$ws.Range("I63").Value = 2500850
$ws.Range("K63").Value = 2500850
$ws.Range("M63").Value = -2500164
# Row 66
$ws.Range("H66").Value = 2500850
$ws.Range("I66").Value = 2500850
$ws.Range("K66").Value = 12504250
$ws.Range("M66").Value = -12500818
# Row 88
$ws.Range("H88").Value = 1006252.8
$ws.Range("I88").Value = 1671302.4
$ws.Range("J88").Value = 8678.5
$ws.Range("K88").Value = 1671302.4
$ws.Range("L88").Value = 8678.5
$ws.Range("M88").Value = -1670896.4
$ws.Range("N88").Value = -9490.5
# Row 91
$ws.Range("H91").Value = 1006252.8
$ws.Range("I91").Value = 1671302.4
$ws.Range("J91").Value = 8678.5
$ws.Range("K91").Value = 1671302.4
$ws.Range("L91").Value = 8678.5
$ws.Range("M91").Value = -1669898.4
$ws.Range("N91").Value = -11486.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 3240
$ws.Range("I86").Value = 3002
$ws.Range("J86").Value = 3537.5
$ws.Range("K86").Value = 3002
$ws.Range("L86").Value = 3537.5
$ws.Range("M86").Value = -1879
$ws.Range("N86").Value = -5783.5
# Row 89
$ws.Range("H89").Value = 3240
$ws.Range("I89").Value = 3002
$ws.Range("J89").Value = 3537.5
$ws.Range("K89").Value = 15010
$ws.Range("L89").Value = 17687.5
$ws.Range("M89").Value = -9394
$ws.Range("N89").Value = -28919.5
# Row 105
$ws.Range("H105").Value = 2940
$ws.Range("I105").Value = 2260
$ws.Range("J105").Value = 3393.3333
$ws.Range("K105").Value = 2260
$ws.Range("L105").Value = 3393.3333
$ws.Range("M105").Value = -513
$ws.Range("N105").Value = -6887.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 19611744
$ws.Range("I62").Value = 3287.7144
$ws.Range("J62").Value = 33337664
$ws.Range("K62").Value = 3287.7144
$ws.Range("L62").Value = 33337664
$ws.Range("M62").Value = -2663.7144
$ws.Range("N62").Value = -33338912
# Row 65
$ws.Range("H65").Value = 19611744
$ws.Range("I65").Value = 3287.7144
$ws.Range("J65").Value = 33337664
$ws.Range("K65").Value = 16438.572
$ws.Range("L65").Value = 166688320
$ws.Range("M65").Value = -13318.572
$ws.Range("N65").Value = -166694560
# Row 132
$ws.Range("H132").Value = 4810950
$ws.Range("I132").Value = 2731.8333
$ws.Range("K132").Value = 8195.499899999999
$ws.Range("M132").Value = -5665.499899999999
# Row 141
$ws.Range("H141").Value = 45386.08
$ws.Range("I141").Value = 17026.285
$ws.Range("J141").Value = 55834.42
$ws.Range("K141").Value = 17026.285
$ws.Range("L141").Value = 55834.42
$ws.Range("M141").Value = -11846.285
$ws.Range("N141").Value = -66194.42

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 1779.81
$ws.Range("I131").Value = 7006.6
$ws.Range("J131").Value = 857.4353
$ws.Range("K131").Value = 21019.8
$ws.Range("L131").Value = 2572.3059
$ws.Range("M131").Value = -15979.8
$ws.Range("N131").Value = -12652.3059

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 52073204
$ws.Range("I70").Value = 62920770
$ws.Range("J70").Value = 4901.8
$ws.Range("K70").Value = 62920770
$ws.Range("L70").Value = 4901.8
$ws.Range("M70").Value = -62920500
$ws.Range("N70").Value = -5441.8
# Row 73
$ws.Range("H73").Value = 52073204
$ws.Range("I73").Value = 62920770
$ws.Range("J73").Value = 4901.8
$ws.Range("K73").Value = 62920770
$ws.Range("L73").Value = 4901.8
$ws.Range("M73").Value = -62919834
$ws.Range("N73").Value = -6773.8
# Row 80
$ws.Range("H80").Value = 3334.2307
$ws.Range("I80").Value = 3606.875
$ws.Range("J80").Value = 2898
$ws.Range("K80").Value = 3606.875
$ws.Range("L80").Value = 2898
$ws.Range("M80").Value = -2608.875
$ws.Range("N80").Value = -4894
# Row 83
$ws.Range("H83").Value = 3334.2307
$ws.Range("I83").Value = 3606.875
$ws.Range("J83").Value = 2898
$ws.Range("K83").Value = 18034.375
$ws.Range("L83").Value = 14490
$ws.Range("M83").Value = -13042.375
$ws.Range("N83").Value = -24474
# Row 97
$ws.Range("H97").Value = 1685.5667
$ws.Range("I97").Value = 1241.762
$ws.Range("K97").Value = 1241.762
$ws.Range("M97").Value = -745.7619999999999
# Row 132
$ws.Range("H132").Value = 2411.6428
$ws.Range("I132").Value = 1877.4
$ws.Range("J132").Value = 3747.25
$ws.Range("K132").Value = 5632.200000000001
$ws.Range("L132").Value = 11241.75
$ws.Range("M132").Value = -3102.200000000001
$ws.Range("N132").Value = -16301.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 5167
$ws.Range("I122").Value = 6886.2856
$ws.Range("J122").Value = 2760
$ws.Range("K122").Value = 20658.8568
$ws.Range("L122").Value = 8280
$ws.Range("M122").Value = -18208.8568
$ws.Range("N122").Value = -13180
